$d = $word.ActiveDocument

# --- Part 1: "9) What is a target?" -> "10) What is a target?" -------------
# Locate the paragraph that starts with "9) What is a target?" (the one with
# the lastRenderedPageBreak, i.e. the *second* "9)" heading in the document).
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("9) What is a target?")) {
        $targetPara = $para
    }
}
if ($targetPara -eq $null) {
    throw "Could not locate the '9) What is a target?' paragraph"
}

$paraRange = $targetPara.Range
# Just the "9) " prefix (3 characters) at the start of the paragraph.
$prefixRange = $d.Range($paraRange.Start, $paraRange.Start + 3)

$renumberXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:lastRenderedPageBreak/><w:t>10</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$prefixRange.InsertXML($renumberXml)

# --- Part 2: append the new Q11 ("What is Virtual DOM?") content -----------
# Insert at the very end of the document body (right before the sectPr),
# i.e. after the last paragraph of Q10 ("A target inside a synthetic ...").
$endPoint = $d.Range($d.Content.End, $d.Content.End)

$newContentXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:ind w:left="-567" w:right="-330"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:ind w:left="-567" w:right="-330"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">11) </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">What is Virtual DOM? </w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="-567" w:right="-330"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>The virtual DOM (VDOM) is a programming concept where an ideal, or “virtual”,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>representation of a UI is kept in memory and synced with the “real” DOM by a library such as</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>React</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202124"/><w:sz w:val="36"/><w:szCs w:val="36"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>DOM. This process is called reconciliation.</w:t></w:r></w:p>

</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$endPoint.InsertXML($newContentXml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
